$d = $word.ActiveDocument

# The paragraph currently reads "Version 1." as:
#   spellStart-proofErr, run("Version"), spellEnd-proofErr,
#   run(" 1."), bookmarkStart("_GoBack"), bookmarkEnd
#
# The target revision keeps the same visible text except "1" -> "2",
# but re-shapes the run/bookmark layout:
#   - "Version" splits into two runs: "Versi" + "on"
#   - " 1." becomes " 2"
#   - a new trailing run holding "." is appended after the bookmark
#
# Locate that paragraph via Find rather than a hard-coded offset.
$search = $d.Content.Duplicate
$found = $search.Find.Execute("Version*.", $true, $false, $true, $false,
                               $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the 'Version 1.' paragraph text"
}

# Re-resolve the owning paragraph straight off the document (rather
# than off the Find duplicate) -- a Range pulled from a duplicated /
# moved Range can report identical Start/End numbers yet still only
# partially honour a later .Delete() (it leaves stray proofErr /
# bookmark markup behind), whereas a Range fetched directly from
# $d.Paragraphs deletes cleanly.
$paraIndex = 1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i).Range
    if ($search.Start -ge $candidate.Start -and $search.Start -lt $candidate.End) {
        $paraIndex = $i
        break
    }
}

# Grab the *whole* paragraph (content + trailing paragraph mark) and
# wipe it -- this clears the text runs plus the proofErr/bookmark
# markup framing them, leaving one clean empty paragraph behind.
$paraRange = $d.Paragraphs.Item($paraIndex).Range
$startPos = $paraRange.Start
$paraRange.Delete()

# Rebuild the paragraph's content from scratch via a raw OOXML
# fragment, so the run/proofErr/bookmark boundaries land exactly
# where the target revision wants them (Find/Replace or Range.Text
# would instead collapse everything back into a single run).
$insertionPoint = $d.Range($startPos, $startPos)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Versi</w:t></w:r><w:r><w:t>on</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 2</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xml)

# InsertXML of a <w:p>-rooted fragment always creates a brand-new
# paragraph just ahead of the (now empty) original paragraph, rather
# than refilling it in place. Delete the paragraph-mark boundary
# between the two so the new content ends up back inside the original
# paragraph element (preserving its w14:paraId / rsid attributes).
$boundary = $d.Range($startPos + 10, $startPos + 11)
$boundary.Delete()
